$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.062.67"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.311.48"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "301.93"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "98.56"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("D10").Value = "35.66"
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "17.91"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "2.671.85"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "2.308.81"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "0.788"
$ws.Range("E17").Value = "  -3.40%  "
$ws.Range("D18").Value = "42.992.02"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  +7.61%  "
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "0.0₃0908"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "68.04"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "241.00"
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("E24").Value = "  -2.94%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "25.05"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").Value = "168.99"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "9.18"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.05"
$ws.Range("E30").Value = "  -8.65%  "
$ws.Range("D31").Value = "33.32"
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "5.20"
$ws.Range("E32").Value = "  +3.92%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").Value = "4.89"
$ws.Range("E33").Value = "  +5.69%  "
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "18.26"
$ws.Range("E35").Value = "  +7.63%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "0.0693"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").Value = "0.102"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").Value = "1.991.74"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").Value = "10.13"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").Value = "2.06"
$ws.Range("E45").Value = "  -10.94%  "
$ws.Range("D46").Value = "17.42"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "2.83"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "76.04"
$ws.Range("E48").Value = "  +8.69%  "
$ws.Range("D49").Value = "54.61"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "2.538.95"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "1.53"
$ws.Range("E51").Value = "  +0.16%  "
